{"js": "// Office.js (Word JavaScript API) implementation of:\n//   \"How/where to deploy an app\"  ->  \"How/where to share an app\"\n// plus relocating the (cursor-position) \"_GoBack\" bookmark so it now sits\n// immediately after the new word \"share\" instead of its old spot later in\n// the document (next to the \"8\" in the schedule table-like list).\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1. Find the word we need to change (\"deploy\" -> \"share\"). Matching the\n//    exact original phrase keeps this targeted to the single occurrence in\n//    the syllabus line \"How/where to deploy an app\".\nconst searchResults = body.search(\"How/where to deploy\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Target text \"How/where to deploy\" was not found.');\n}\n\nconst targetRange = searchResults.items[0];\ntargetRange.insertText(\"How/where to share\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Move the \"_GoBack\" bookmark: remove it from wherever it currently is\n//    and re-insert it (collapsed, i.e. zero-length) right after the word\n//    \"share\" we just inserted.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst shareResults = body.search(\"share\", { matchCase: true });\nshareResults.load(\"items\");\nawait context.sync();\n\nif (shareResults.items.length === 0) {\n  throw new Error('Could not locate \"share\" after insertion.');\n}\n\nconst shareRange = shareResults.items[shareResults.items.length - 1];\nconst afterShare = shareRange.getRange(Word.RangeLocation.after);\nafterShare.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop implementation of:\n#   \"How/where to deploy an app\"  ->  \"How/where to share an app\"\n# plus relocating the (cursor-position) \"_GoBack\" bookmark so it now sits\n# immediately after the new word \"share\" instead of its old spot later in\n# the document (next to the \"8\" in the schedule table-like list).\n\n$d = $word.ActiveDocument\n\n# 1. Locate the exact phrase and change \"deploy\" to \"share\" in place.\n$target = $d.Content\n$target.Find.ClearFormatting()\n$found = $target.Find.Execute(\"How/where to deploy\")\nif (-not $found) {\n    throw 'Target text \"How/where to deploy\" was not found.'\n}\n$target.Text = \"How/where to share\"\n\n# 2. Move the \"_GoBack\" bookmark: remove whatever copy currently exists and\n#    add a new, collapsed (zero-length) one right after the word \"share\".\n$shareRange = $d.Content\n$shareRange.Find.ClearFormatting()\n$shareRange.Find.Execute(\"share\") | Out-Null\n$shareRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $shareRange) | Out-Null\n"}
